# Update column G ("K") values on the active sheet per the new scraped data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 5
    3  = 2
    4  = 7
    5  = 7
    6  = 4
    7  = 1
    8  = 3
    9  = 6
    10 = 1
    11 = 3
    12 = 5
    13 = 6
    14 = 5
    15 = 5
    16 = 4
    17 = 6
    18 = 2
    19 = 13
    20 = 8
    21 = 9
    22 = 6
    23 = 7
    24 = 4
    25 = 8
    26 = 8
    27 = 6
    28 = 4
    29 = 3
    30 = 3
    31 = 5
    32 = 7
    33 = 1
    34 = 1
    35 = 5
    36 = 3
    37 = 3
    38 = 1
}

foreach ($row in $newValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $newValues[$row]
}
